$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.147.24'
$ws.Range("E2").Value = '  -2.16%  '
$ws.Range("D3").Value = '1.577.53'
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '209.22'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("E8").Value = '  -0.61%  '
$ws.Range("E9").Value = '  -1.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.53'
$ws.Range("E10").Value = '  -0.71%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0843'
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").Value = '1.799.03'
$ws.Range("E12").Value = '  -1.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.05'
$ws.Range("E13").Value = '  -0.06%  '
$ws.Range("D14").Value = '1.536.44'
$ws.Range("E14").Value = '  -3.95%  '
$ws.Range("E15").Value = '  -1.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.43'
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").Value = '26.153.44'
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.26'
$ws.Range("E19").Value = '  +1.68%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '207.77'
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("E22").Value = '  -1.00%  '
$ws.Range("E23").Value = '  -2.68%  '
$ws.Range("E24").Value = '  -1.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.86'
$ws.Range("E25").Value = '  +0.20%  '
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.98'
$ws.Range("E27").Value = '  -1.37%  '
$ws.Range("E28").Value = '  -1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.23'
$ws.Range("E29").Value = '  -0.70%  '
$ws.Range("E30").Value = '  -0.78%  '
$ws.Range("E31").Value = '  -1.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.20'
$ws.Range("E32").Value = '  -1.69%  '
$ws.Range("E33").Value = '  +0.28%  '
$ws.Range("D34").Value = '1.277.72'
$ws.Range("E34").Value = '  -0.78%  '
$ws.Range("E35").Value = '  -0.57%  '
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0165'
$ws.Range("E38").Value = '  -2.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.09'
$ws.Range("E39").Value = '  -6.14%  '
$ws.Range("E40").Value = '  -1.52%  '
$ws.Range("E41").Value = '  +3.00%  '
$ws.Range("E42").Value = '  -2.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.46'
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("E44").Value = '  -2.72%  '
$ws.Range("D45").Value = '1.712.48'
$ws.Range("E45").Value = '  -1.45%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.94'
$ws.Range("E46").Value = '  -1.68%  '
$ws.Range("E47").Value = '  -0.52%  '
$ws.Range("E48").Value = '  -1.73%  '
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0506'
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.74'
$ws.Range("E51").Value = '  +10.86%  '
